$wb = $excel.ActiveWorkbook

# --- Sheet "Games": append the newly-played game (row 41) ---
$games = $wb.Worksheets.Item("Games")

$row = 41
$games.Range("A$row").Value = 40
$games.Range("B$row").NumberFormat = "YYYY-MM-DD"
$games.Range("B$row").Value = 45307
$games.Range("C$row").Value = 1
$games.Range("D$row").Value = 128
$games.Range("E$row").Value = 91.90000000000001
$games.Range("F$row").Value = 0.655
$games.Range("G$row").Value = 8.5
$games.Range("H$row").Value = 31.6
$games.Range("I$row").Value = 0.161
$games.Range("J$row").Value = 139.3
$games.Range("K$row").Value = "OKC"
$games.Range("L$row").Value = 117
$games.Range("M$row").Value = 0.602
$games.Range("N$row").Value = 11.6
$games.Range("O$row").Value = 25
$games.Range("P$row").Value = 0.205
$games.Range("Q$row").Value = 127.4
$games.Range("R$row").Value = 1
$games.Range("S$row").Value = 1

# --- Sheet "Next": the game that was just played (old row 2, OKC on 1/16)
#     drops off the upcoming-schedule list; everything below shifts up ---
$next = $wb.Worksheets.Item("Next")
$next.Rows.Item(2).Delete()
